$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 16 switches to a different (built-in) table style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{820FE5FD-D9EC-4A69-A41C-2FC222666D55}")

# ---------------------------------------------------------------------------
# 2) The deck's theme (applied through the slide master) switches from the
#    custom "Integral" theme to the stock "Office Theme" palette.
# ---------------------------------------------------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

function BGR($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$themeColors.Colors(1).RGB  = BGR 0x00 0x00 0x00   # dk1
$themeColors.Colors(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1
$themeColors.Colors(3).RGB  = BGR 0x44 0x54 0x6A   # dk2
$themeColors.Colors(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2
$themeColors.Colors(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1
$themeColors.Colors(6).RGB  = BGR 0xED 0x7D 0x31   # accent2
$themeColors.Colors(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3
$themeColors.Colors(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4
$themeColors.Colors(9).RGB  = BGR 0x44 0x72 0xC4   # accent5
$themeColors.Colors(10).RGB = BGR 0x70 0xAD 0x47   # accent6
$themeColors.Colors(11).RGB = BGR 0x05 0x63 0xC1   # hlink
$themeColors.Colors(12).RGB = BGR 0x95 0x4F 0x72   # folHlink
